$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "코사인 유사도(cosine similarity)로 과거 주가의 유사 패턴을 찾아 미래 예측하기"
$ws.Range("E4").Value = "https://teddylee777.github.io/pandas/cos-sim-stock"

$ws.Range("D12").Value = "Yes24 2021년 올해의 책 후보 도전에 응원 부탁드립니다."
$ws.Range("E12").Value = "https://tensorflow.blog/2021/10/30/yes24-2021%eb%85%84-%ec%98%ac%ed%95%b4%ec%9d%98-%ec%b1%85-%ed%9b%84%eb%b3%b4-%eb%8f%84%ec%a0%84%ec%97%90-%ec%9d%91%ec%9b%90-%eb%b6%80%ed%83%81%eb%93%9c%eb%a6%bd%eb%8b%88%eb%8b%a4/"

$ws.Range("D20").Value = "[머신러닝 기초] 지도학습 - 선형 회귀(Regression) 분석"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/585"

$ws.Range("D28").Value = "강화학습 실습 2편 : OpenAI GYM 기반 환경(env) 코드 분석"

$ws.Range("D44").Value = "Qualcomm Ventures Portfolio (3) - Azion"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/104"

$ws.Range("D46").Value = "[Bioinformatics] 2021년 11월,  유전체 정보분석 전문가 심화과정"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/419"
